# Balmorel plotting tool: expand manual_colors_input workbook to cover
# both Electricity and Heat balances ("HOURLY BALANCE" plotting improvement).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename + rebuild the three existing sheets (Production/Demand/Price ->
#    Electricity_Production/Electricity_Demand/Electricity_Price) with their
#    refreshed category -> color tables.
# ---------------------------------------------------------------------------

$wsElecProd = $wb.Worksheets.Item(1)
$wsElecProd.Name = "Electricity_Production"
$wsElecProd.Range("A1:B100").ClearContents()
$wsElecProd.Range("A1").Value = "Value"
$wsElecProd.Range("B1").Value = "Color"
$wsElecProd.Range("A2").Value = "THERMAL"
$wsElecProd.Range("B2").Value = "violet"
$wsElecProd.Range("A3").Value = "ELSTO"
$wsElecProd.Range("B3").Value = "grey"
$wsElecProd.Range("A4").Value = "HYDRO"
$wsElecProd.Range("B4").Value = "dark blue"
$wsElecProd.Range("A5").Value = "WIND"
$wsElecProd.Range("B5").Value = "lightgreen"
$wsElecProd.Range("A6").Value = "SOLARPV"
$wsElecProd.Range("B6").Value = "yellow"
$wsElecProd.Range("A7").Value = "FUELCELL"
$wsElecProd.Range("B7").Value = "violet"

$wsElecDem = $wb.Worksheets.Item(2)
$wsElecDem.Name = "Electricity_Demand"
$wsElecDem.Range("A1:B100").ClearContents()
$wsElecDem.Range("A1").Value = "Value"
$wsElecDem.Range("B1").Value = "Color"
$wsElecDem.Range("A2").Value = "DEMAND_EXO"
$wsElecDem.Range("B2").Value = "blue"
$wsElecDem.Range("A3").Value = "DEMAND_LOSS"
$wsElecDem.Range("B3").Value = "red"
$wsElecDem.Range("A4").Value = "DEMAND_ELSTO"
$wsElecDem.Range("B4").Value = "orange"
$wsElecDem.Range("A5").Value = "DEMAND_P2H"
$wsElecDem.Range("B5").Value = "dark grey"
$wsElecDem.Range("A6").Value = "DEMAND_EV"
$wsElecDem.Range("B6").Value = "green"
$wsElecDem.Range("A7").Value = "DEMAND_P2G"
$wsElecDem.Range("B7").Value = "dark blue"
$wsElecDem.Range("A8").Value = "DEMAND_CCS"
$wsElecDem.Range("B8").Value = "pink"

$wsElecPrice = $wb.Worksheets.Item(3)
$wsElecPrice.Name = "Electricity_Price"
$wsElecPrice.Range("A1:B100").ClearContents()
$wsElecPrice.Range("A1").Value = "Value"
$wsElecPrice.Range("B1").Value = "Color"
$wsElecPrice.Range("A2").Value = "Price"
$wsElecPrice.Range("B2").Value = "Black"

# ---------------------------------------------------------------------------
# 2) Add the three new Heat_* sheets at the end of the workbook, mirroring
#    the structure of their Electricity_* counterparts.
# ---------------------------------------------------------------------------

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHeatProd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$wsHeatProd.Name = "Heat_Production"
$wsHeatProd.Range("A1").Value = "Value"
$wsHeatProd.Range("B1").Value = "Color"
$wsHeatProd.Range("A2").Value = "BOILERS"
$wsHeatProd.Range("B2").Value = "violet"
$wsHeatProd.Range("A3").Value = "CHP"
$wsHeatProd.Range("B3").Value = "grey"
$wsHeatProd.Range("A4").Value = "SOLARHEATING"
$wsHeatProd.Range("B4").Value = "dark blue"
$wsHeatProd.Range("A5").Value = "HEATSTO"
$wsHeatProd.Range("B5").Value = "lightgreen"
$wsHeatProd.Range("A6").Value = "P2H"
$wsHeatProd.Range("B6").Value = "yellow"
$wsHeatProd.Range("A7").Value = "FUELCELL"
$wsHeatProd.Range("B7").Value = "violet"
$wsHeatProd.Activate() | Out-Null
$excel.ActiveWindow.DisplayGridlines = $false
$wsHeatProd.Range("A2:A7").Select() | Out-Null

$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHeatDem = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet2)
$wsHeatDem.Name = "Heat_Demand"
$wsHeatDem.Range("A1").Value = "Value"
$wsHeatDem.Range("B1").Value = "Color"
$wsHeatDem.Range("A2").Value = "DEMAND_EXO"
$wsHeatDem.Range("B2").Value = "blue"
$wsHeatDem.Range("A3").Value = "DEMAND_HEATSTO"
$wsHeatDem.Range("B3").Value = "red"
$wsHeatDem.Range("A4").Value = "DEMAND_LOSS"
$wsHeatDem.Range("B4").Value = "pink"
$wsHeatDem.Range("A5").Value = "DEMAND_P2G"
$wsHeatDem.Range("B5").Value = "green"
$wsHeatDem.Activate() | Out-Null
$excel.ActiveWindow.DisplayGridlines = $false
$wsHeatDem.Range("A3").Select() | Out-Null

$afterSheet3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHeatPrice = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet3)
$wsHeatPrice.Name = "Heat_Price"
$wsHeatPrice.Range("A1").Value = "Value"
$wsHeatPrice.Range("B1").Value = "Color"
$wsHeatPrice.Range("A2").Value = "Price"
$wsHeatPrice.Range("B2").Value = "Black"
$wsHeatPrice.Activate() | Out-Null
$excel.ActiveWindow.DisplayGridlines = $false
$wsHeatPrice.Range("A1:B2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-select the full data range on the three original sheets (mirrors the
#    "select-all-before-save" state captured in the saved workbook) and make
#    Heat_Demand the active/visible tab, matching the final saved view.
# ---------------------------------------------------------------------------

$wsElecProd.Activate() | Out-Null
$wsElecProd.Range("A1:B7").Select() | Out-Null

$wsElecDem.Activate() | Out-Null
$wsElecDem.Range("A1:B8").Select() | Out-Null

$wsElecPrice.Activate() | Out-Null
$wsElecPrice.Range("A1:B2").Select() | Out-Null

$wsHeatDem.Activate() | Out-Null
$wsHeatDem.Range("A3").Select() | Out-Null
